# Adding solution to the "Longest repeating character replacement" problem.
#
# This mirrors the author's edit: a new row (15) is appended to the tracker
# sheet with the problem name in column B (same style/pattern as the other
# rows in that column), the view is scrolled/zoomed to the new row, and the
# column widths for A and B are nudged slightly narrower/wider.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New row: B15 gets the new problem name (same "s=1" wrap-text style the
#    rest of column B already uses - Excel applies the column's style
#    automatically for a cell with no prior formatting).
$ws.Range("B15").Value = "Longest repeating character replacement"

# 2. Column width tweaks (col A: 21.66 -> 16.66 "chars", col B: 38.66 -> 44.16
#    "chars"). Feed ColumnWidth the same character-width numbers Excel itself
#    would report for those target widths.
$ws.Columns.Item(1).ColumnWidth = 15.833333333333334
$ws.Columns.Item(2).ColumnWidth = 43.333333333333336

# 3. Update the view: zoom in to 262% and scroll/select down to the newly
#    added row so B15 is the active cell.
$excel.ActiveWindow.Zoom = 262
$ws.Range("B15").Select()
